$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "WellPad_Ac" -> "WellPad_Active" for the well-pad rows (A9:A13)
$ws.Range("A9:A13").Value = "WellPad_Active"

# Move the selection, matching the saved cursor position in the target file
$ws.Range("A5").Select()
